$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.439248
$ws.Range("H2").Value = 13.317744
$ws.Range("I2").Value = 0.3193177756555054
$ws.Range("J2").Value = 0.3193177756555054
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 104.9251675201547
$ws.Range("R2").Value = 944.326507681392
$ws.Range("S2").Value = 0.02180252146240675
$ws.Range("T2").Value = 0.02180252146240675

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.439248
$ws.Range("H3").Value = 13.317744
$ws.Range("I3").Value = 0.3193177756555054
$ws.Range("J3").Value = 0.3193177756555054
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 804.7841286497174
$ws.Range("R3").Value = 7243.057157847456
$ws.Range("S3").Value = 0.1672270214304817
$ws.Range("T3").Value = 0.1672270214304816

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.439248
$ws.Range("H4").Value = 13.317744
$ws.Range("I4").Value = 0.3193177756555054
$ws.Range("J4").Value = 0.3193177756555054
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 493.2028781091573
$ws.Range("R4").Value = 4438.825902982416
$ws.Range("S4").Value = 0.1024831943511567
$ws.Range("T4").Value = 0.1024831943511567

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.439248
$ws.Range("H5").Value = 13.317744
$ws.Range("I5").Value = 0.3193177756555054
$ws.Range("J5").Value = 0.3193177756555054
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 133.8124270744213
$ws.Range("R5").Value = 1204.311843669792
$ws.Range("S5").Value = 0.02780503841146032
$ws.Range("T5").Value = 0.02780503841146032

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.893859666666666
$ws.Range("H6").Value = 11.681579
$ws.Range("I6").Value = 0.2800876651799331
$ws.Range("J6").Value = 0.2800876651799331
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 92.03447922372743
$ws.Range("R6").Value = 828.310313013547
$ws.Range("S6").Value = 0.01912395048758258
$ws.Range("T6").Value = 0.01912395048758258

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.893859666666666
$ws.Range("H7").Value = 11.681579
$ws.Range("I7").Value = 0.2800876651799331
$ws.Range("J7").Value = 0.2800876651799331
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 705.9115550477495
$ws.Range("R7").Value = 6353.203995429746
$ws.Range("S7").Value = 0.1466821754326306
$ws.Range("T7").Value = 0.1466821754326306

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.893859666666666
$ws.Range("H8").Value = 11.681579
$ws.Range("I8").Value = 0.2800876651799331
$ws.Range("J8").Value = 0.2800876651799331
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 432.6099363119978
$ws.Range("R8").Value = 3893.48942680798
$ws.Range("S8").Value = 0.08989251715496185
$ws.Range("T8").Value = 0.08989251715496187

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.893859666666666
$ws.Range("H9").Value = 11.681579
$ws.Range("I9").Value = 0.2800876651799331
$ws.Range("J9").Value = 0.2800876651799331
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 117.3727650908135
$ws.Range("R9").Value = 1056.354885817322
$ws.Range("S9").Value = 0.02438902210475799
$ws.Range("T9").Value = 0.024389022104758

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.735205333333333
$ws.Range("H10").Value = 5.205616
$ws.Range("I10").Value = 0.1248143621049263
$ws.Range("J10").Value = 0.1248143621049263
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 41.01296216878755
$ws.Range("R10").Value = 369.116659519088
$ws.Range("S10").Value = 0.008522130667555104
$ws.Range("T10").Value = 0.008522130667555104

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.735205333333333
$ws.Range("H11").Value = 5.205616
$ws.Range("I11").Value = 0.1248143621049263
$ws.Range("J11").Value = 0.1248143621049263
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 314.5725835130204
$ws.Range("R11").Value = 2831.153251617184
$ws.Range("S11").Value = 0.06536539960453198
$ws.Range("T11").Value = 0.06536539960453197

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.735205333333333
$ws.Range("H12").Value = 5.205616
$ws.Range("I12").Value = 0.1248143621049263
$ws.Range("J12").Value = 0.1248143621049263
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 192.7822605338471
$ws.Range("R12").Value = 1735.040344804624
$ws.Range("S12").Value = 0.04005844805587874
$ws.Range("T12").Value = 0.04005844805587875

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.735205333333333
$ws.Range("H13").Value = 5.205616
$ws.Range("I13").Value = 0.1248143621049263
$ws.Range("J13").Value = 0.1248143621049263
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 52.30436261407644
$ws.Range("R13").Value = 470.7392635266879
$ws.Range("S13").Value = 0.01086838377696045
$ws.Range("T13").Value = 0.01086838377696045

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.833976
$ws.Range("H14").Value = 11.501928
$ws.Range("I14").Value = 0.2757801970596353
$ws.Range("J14").Value = 0.2757801970596353
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 90.61908099485599
$ws.Range("R14").Value = 815.571728953704
$ws.Range("S14").Value = 0.01882984325866732
$ws.Range("T14").Value = 0.01882984325866732

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.833976
$ws.Range("H15").Value = 11.501928
$ws.Range("I15").Value = 0.2757801970596353
$ws.Range("J15").Value = 0.2757801970596353
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 695.055341450608
$ws.Range("R15").Value = 6255.498073055472
$ws.Range("S15").Value = 0.1444263502998599
$ws.Range("T15").Value = 0.1444263502998598

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.833976
$ws.Range("H16").Value = 11.501928
$ws.Range("I16").Value = 0.2757801970596353
$ws.Range("J16").Value = 0.2757801970596353
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 425.956828228888
$ws.Range("R16").Value = 3833.611454059992
$ws.Range("S16").Value = 0.08851006016011502
$ws.Range("T16").Value = 0.08851006016011502

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.833976
$ws.Range("H17").Value = 11.501928
$ws.Range("I17").Value = 0.2757801970596353
$ws.Range("J17").Value = 0.2757801970596353
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 115.567689371056
$ws.Range("R17").Value = 1040.109204339504
$ws.Range("S17").Value = 0.02401394334099311
$ws.Range("T17").Value = 0.02401394334099311
